# Updates network stack graphic to say "OSI Model" instead of "OSI".
#
# The text lives in a standalone textbox on slide 1 (a single run, no
# other runs share the paragraph), so we find it by its exact current
# text and update the TextRange in place. Setting .Text on a TextRange
# that already contains exactly one uniformly-formatted run preserves
# that run's formatting (font, size, bold, etc.), so only the characters
# change - matching the diff, which leaves <a:rPr> untouched and only
# rewrites the <a:t> content.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -eq "OSI") {
            $tr.Text = "OSI Model"
        }
    }
}
